# Update the "Restricciones_del_follower" sheet data rows (2-6)
$wb = $excel.ActiveWorkbook

$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")

$wsFollower.Range("A2").Value = "7.35 - y_1"
$wsFollower.Range("B2").Value = "-7.35"
$wsFollower.Range("C2").Value = "J_0_L0_v"
$wsFollower.Range("D2").Value = "0.03"
$wsFollower.Range("E2").Value = "-2.8000000000000003"
$wsFollower.Range("F2").Value = "-1.3"

$wsFollower.Range("A3").Value = "-7.35 + y_1"
$wsFollower.Range("B3").Value = "3.3499999999999996"
$wsFollower.Range("C3").Value = "J_0_L0_v"
$wsFollower.Range("D3").Value = "0.85"
$wsFollower.Range("E3").Value = "-2.1"
$wsFollower.Range("F3").Value = "-7.199999999999999"

$wsFollower.Range("A4").Value = "-3.7499999999999964 - 2x + y_1 + 4y_2"
$wsFollower.Range("B4").Value = "-12.250000000000004"
$wsFollower.Range("C4").Value = "J_0_LP_v"
$wsFollower.Range("D4").Value = "0.73"
$wsFollower.Range("E4").Value = "-8.8"
$wsFollower.Range("F4").Value = "-6.4"

$wsFollower.Range("A5").Value = "-67.22 + 8x + y_1"
$wsFollower.Range("B5").Value = "18.550000000000004"
$wsFollower.Range("C5").Value = "J_Ne_L0_v"
$wsFollower.Range("D5").Value = "0.02"
$wsFollower.Range("E5").Value = "1.6"
$wsFollower.Range("F5").Value = "6.800000000000001"

$wsFollower.Range("A6").Value = "-5.5 - 2x - 2y_1"
$wsFollower.Range("B6").Value = "-17.5"
$wsFollower.Range("C6").Value = "J_Ne_L0_v"
$wsFollower.Range("D6").Value = "0.72"
$wsFollower.Range("E6").Value = "-0.0"
$wsFollower.Range("F6").Value = "-5.5"

# Update the "Punto_modificado" sheet (modified point x, y_1, y_2)
$wsPunto = $wb.Worksheets.Item("Punto_modificado")
$wsPunto.Range("A2").Value = "7.4"
$wsPunto.Range("B2").Value = "7.35"
$wsPunto.Range("C2").Value = "2.8"

# Update the "Vector_bf" sheet (lowercase bf) - sheet names "Vector_bf" and
# "Vector_BF" differ only by case, and worksheet lookup by name is
# case-insensitive, so address these two specifically by their (1-based)
# tab position instead of by name to avoid ambiguity.
$wsVecLower = $wb.Worksheets.Item(5)   # Vector_bf
$wsVecLower.Range("A2").Value = "0.8700000000000001"
$wsVecLower.Range("A3").Value = "-2.92"

# Update the "Vector_BF" sheet (uppercase BF)
$wsVecUpper = $wb.Worksheets.Item(6)   # Vector_BF
$wsVecUpper.Range("A2").Value = "-29.400000000000002"
$wsVecUpper.Range("A3").Value = "9.500000000000002"
$wsVecUpper.Range("A4").Value = "33.2"

Write-Host "edit applied"
